$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 55 (shifts existing rows 55-89 down to 57-91)
$ws.Rows.Item(55).Resize(2).Insert()

# Populate the two new rows with the new ticker symbols (only column A is filled)
$ws.Range("A55").Value = "MNST.US"
$ws.Range("A56").Value = "MSFT.US"
